$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H40").Value = 3209.9
$ws.Range("I40").Value = 1871.4286
$ws.Range("J40").Value = 6333
$ws.Range("K40").Value = 1871.4286
$ws.Range("L40").Value = 6333
$ws.Range("M40").Value = -1696.4286
$ws.Range("N40").Value = -6683
$ws.Range("H64").Value = 3773.3
$ws.Range("I64").Value = 3558.5862
$ws.Range("K64").Value = 3558.5862
$ws.Range("M64").Value = -3310.5862
$ws.Range("H67").Value = 3773.3
$ws.Range("I67").Value = 3558.5862
$ws.Range("K67").Value = 3558.5862
$ws.Range("M67").Value = -2700.5862
$ws.Range("H74").Value = 4299.6924
$ws.Range("J74").Value = 5142.857
$ws.Range("L74").Value = 5142.857
$ws.Range("N74").Value = -7014.857
$ws.Range("H77").Value = 4299.6924
$ws.Range("J77").Value = 5142.857
$ws.Range("L77").Value = 25714.285
$ws.Range("N77").Value = -35074.285
$ws.Range("H80").Value = 1732.2285
$ws.Range("I80").Value = 796.8461
$ws.Range("J80").Value = 2284.9546
$ws.Range("K80").Value = 2390.5383
$ws.Range("L80").Value = 6854.8638
$ws.Range("M80").Value = -1392.5383
$ws.Range("N80").Value = -8850.863799999999
$ws.Range("H83").Value = 1732.2285
$ws.Range("I83").Value = 796.8461
$ws.Range("J83").Value = 2284.9546
$ws.Range("K83").Value = 7171.6149
$ws.Range("L83").Value = 20564.5914
$ws.Range("M83").Value = -2179.6149
$ws.Range("N83").Value = -30548.5914
$ws.Range("H111").Value = 750
$ws.Range("I111").Value = 500
$ws.Range("J111").Value = 1000
$ws.Range("K111").Value = 1500
$ws.Range("L111").Value = 3000
$ws.Range("M111").Value = 1567
$ws.Range("N111").Value = -9134
$ws.Range("H125").Value = 1537.125
$ws.Range("I125").Value = 799.55554
$ws.Range("K125").Value = 7195.99986
$ws.Range("M125").Value = -4735.99986
$ws.Range("H132").Value = 2191.7778
$ws.Range("I132").Value = 2203
$ws.Range("J132").Value = 1900
$ws.Range("K132").Value = 6609
$ws.Range("L132").Value = 5700
$ws.Range("M132").Value = -4079
$ws.Range("N132").Value = -10760
$ws.Range("H134").Value = 100413.89
$ws.Range("J134").Value = 90528.375
$ws.Range("L134").Value = 90528.375
$ws.Range("N134").Value = -100668.375
$ws.Range("H137").Value = 1137.238
$ws.Range("I137").Value = 1100.2307
$ws.Range("J137").Value = 1197.375
$ws.Range("K137").Value = 3300.6921
$ws.Range("L137").Value = 3592.125
$ws.Range("M137").Value = -750.6921000000002
$ws.Range("N137").Value = -8692.125
$ws.Range("H138").Value = 2575.5806
$ws.Range("I138").Value = 2507.5
$ws.Range("J138").Value = 2608
$ws.Range("K138").Value = 7522.5
$ws.Range("L138").Value = 7824
$ws.Range("M138").Value = -2382.5
$ws.Range("N138").Value = -18104
$ws.Range("H141").Value = 2389.2727
$ws.Range("I141").Value = 2519.2
$ws.Range("J141").Value = 1090
$ws.Range("K141").Value = 7557.599999999999
$ws.Range("L141").Value = 3270
$ws.Range("M141").Value = -2377.599999999999
$ws.Range("N141").Value = -13630

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H4").Value = 286031.16
$ws.Range("I4").Value = 286031.16
$ws.Range("K4").Value = 286031.16
$ws.Range("M4").Value = -285915.16
$ws.Range("H32").Value = 21936.688
$ws.Range("I32").Value = 6915.6665
$ws.Range("K32").Value = 6915.6665
$ws.Range("M32").Value = -6628.6665
$ws.Range("H61").Value = 1847.9565
$ws.Range("I61").Value = 1800.1111
$ws.Range("K61").Value = 1800.1111
$ws.Range("M61").Value = -1588.1111
$ws.Range("H136").Value = 1847.9565
$ws.Range("I136").Value = 1800.1111
$ws.Range("K136").Value = 5400.3333
$ws.Range("M136").Value = -2850.3333

# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 10050.973
$ws.Range("I20").Value = 10112.27
$ws.Range("J20").Value = 9891.6
$ws.Range("K20").Value = 10112.27
$ws.Range("L20").Value = 9891.6
$ws.Range("M20").Value = -9865.27
$ws.Range("N20").Value = -10385.6
$ws.Range("H99").Value = 2381.9
$ws.Range("I99").Value = 1688.4286
$ws.Range("K99").Value = 1688.4286
$ws.Range("M99").Value = -190.4286
$ws.Range("H134").Value = 1408.069
$ws.Range("I134").Value = 1171.3043
$ws.Range("K134").Value = 3513.9129
$ws.Range("M134").Value = -978.9129000000003

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H141").Value = 270211.9
$ws.Range("J141").Value = 270211.9
$ws.Range("L141").Value = 270211.9
$ws.Range("N141").Value = -280571.9

# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H14").Value = 1529.579
$ws.Range("I14").Value = 1529.579
$ws.Range("K14").Value = 4588.737
$ws.Range("M14").Value = -4415.737
$ws.Range("H44").Value = 0
$ws.Range("I44").Value = 0
$ws.Range("J44").Value = 0
$ws.Range("K44").Value = 0
$ws.Range("L44").Value = 0
$ws.Range("M44").ClearContents()
$ws.Range("N44").ClearContents()
$ws.Range("H63").Value = 5768
$ws.Range("I63").Value = 1750
$ws.Range("J63").Value = 7777
$ws.Range("K63").Value = 5250
$ws.Range("L63").Value = 23331
$ws.Range("M63").Value = -4501
$ws.Range("N63").Value = -24829
$ws.Range("H64").Value = 2120.7273
$ws.Range("I64").Value = 1913.1666
$ws.Range("J64").Value = 2369.8
$ws.Range("K64").Value = 5739.4998
$ws.Range("L64").Value = 7109.400000000001
$ws.Range("M64").Value = -5469.4998
$ws.Range("N64").Value = -7649.400000000001
$ws.Range("H66").Value = 5768
$ws.Range("I66").Value = 1750
$ws.Range("J66").Value = 7777
$ws.Range("K66").Value = 15750
$ws.Range("L66").Value = 69993
$ws.Range("M66").Value = -12006
$ws.Range("N66").Value = -77481
$ws.Range("H67").Value = 2120.7273
$ws.Range("I67").Value = 1913.1666
$ws.Range("J67").Value = 2369.8
$ws.Range("K67").Value = 5739.4998
$ws.Range("L67").Value = 7109.400000000001
$ws.Range("M67").Value = -4803.4998
$ws.Range("N67").Value = -8981.400000000001
$ws.Range("H82").Value = 8750
$ws.Range("J82").Value = 8750
$ws.Range("L82").Value = 26250
$ws.Range("N82").Value = -27062
$ws.Range("H85").Value = 8750
$ws.Range("J85").Value = 8750
$ws.Range("L85").Value = 26250
$ws.Range("N85").Value = -29058
$ws.Range("H113").Value = 2539.4
$ws.Range("J113").Value = 3666.3333
$ws.Range("L113").Value = 10998.9999
$ws.Range("N113").Value = -15338.9999
$ws.Range("H121").Value = 17598022
$ws.Range("J121").Value = 85689.75
$ws.Range("L121").Value = 257069.25
$ws.Range("N121").Value = -259689.25
$ws.Range("H131").Value = 1915.75
$ws.Range("I131").Value = 2176.6667
$ws.Range("K131").Value = 6530.000100000001
$ws.Range("M131").Value = -1490.000100000001

# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H20").Value = 40957.31
$ws.Range("J20").Value = 44767.273
$ws.Range("L20").Value = 44767.273
$ws.Range("N20").Value = -45257.273
$ws.Range("H24").Value = 30519.375
$ws.Range("J24").Value = 31820.666
$ws.Range("L24").Value = 31820.666
$ws.Range("N24").Value = -32166.666
$ws.Range("H70").Value = 11665.5
$ws.Range("I70").Value = 11665.5
$ws.Range("K70").Value = 11665.5
$ws.Range("M70").Value = -11395.5
$ws.Range("H73").Value = 11665.5
$ws.Range("I73").Value = 11665.5
$ws.Range("K73").Value = 11665.5
$ws.Range("M73").Value = -10729.5
$ws.Range("H132").Value = 3884.75
$ws.Range("I132").Value = 3997.0667
$ws.Range("K132").Value = 11991.2001
$ws.Range("M132").Value = -9461.2001

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H16").Value = 84604.664
$ws.Range("I16").Value = 1265.6
$ws.Range("J16").Value = 501300
$ws.Range("K16").Value = 1265.6
$ws.Range("L16").Value = 501300
$ws.Range("M16").Value = -1095.6
$ws.Range("N16").Value = -501640
$ws.Range("H61").Value = 80759.88
$ws.Range("I61").Value = 92871.73
$ws.Range("K61").Value = 92871.73
$ws.Range("M61").Value = -92669.73
$ws.Range("H113").Value = 80759.88
$ws.Range("I113").Value = 92871.73
$ws.Range("K113").Value = 92871.73
$ws.Range("M113").Value = -90701.73
$ws.Range("H136").Value = 4490
$ws.Range("I136").Value = 4221.6
$ws.Range("J136").Value = 5026.8
$ws.Range("K136").Value = 12664.8
$ws.Range("L136").Value = 15080.4
$ws.Range("M136").Value = -10114.8
$ws.Range("N136").Value = -20180.4

# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H132").Value = 1344.1136
$ws.Range("I132").Value = 918.5172
$ws.Range("J132").Value = 2166.9333
$ws.Range("K132").Value = 2755.5516
$ws.Range("L132").Value = 6500.7999
$ws.Range("M132").Value = -225.5515999999998
$ws.Range("N132").Value = -11560.7999
$ws.Range("H136").Value = 2728.342
$ws.Range("I136").Value = 2583.7036
$ws.Range("J136").Value = 3083.3635
$ws.Range("K136").Value = 7751.110799999999
$ws.Range("L136").Value = 9250.0905
$ws.Range("M136").Value = -5201.110799999999
$ws.Range("N136").Value = -14350.0905
